$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 115, pushing the old "footer" row (note about the
# 4/8 merge of the two consultation windows) down to row 116, carrying its
# formatting with it.
$ws.Rows.Item(115).Insert()

# Fill the newly-freed row 115 with the new day's data (2020-05-19).
$ws.Cells.Item(115, 1).Value = 43970
$ws.Cells.Item(115, 2).Value = 206
$ws.Cells.Item(115, 3).Value = 38171
$ws.Cells.Item(115, 4).Value = 40
$ws.Cells.Item(115, 5).Value = 7682

# Extend the print area to include the (now pushed-down) footer row.
$pa = $wb.Names.Item(1)
$pa.RefersTo = '=相談件数!$A$1:$E$116'

# Move the selection in the frozen (bottom-right) pane to match the
# author's final cursor position.
$ws.Range("D113").Select()
